$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Row 11 ----
$ws.Range("A11").Value = "11. Fast & Slow Pointers"
$ws.Range("B11").Value = "middleNode"
$ws.Range("C11").Value = "Go throw the linked list with fast & slow pointer. The fast pointer is going two nodes at one jump and the slow one node at the time." + [char]10 + "Meaning when the fast pointer gets to the end the slow pointer will be doing half of the way hence middle."

# ---- Row 12 ----
$ws.Range("A12").Value = "12. Fast &Slow Pointers"
$ws.Range("B12").Value = "Palindrome Linked list"
$ws.Range("C12").Value = "find the middle point by question (11) and save the values in array. Compare the rest of the linked list to the array's values in reverse"

# ---- Row 13 ----
$ws.Range("A13").Value = "13. Fast & Slow Pointers"
$ws.Range("B13").Value = "remove elements"
$ws.Range("C13").Value = "If we have a head then call stack to the recursive function. After getting to the end check in reverse if the current value need to be removed if it is remove it." + [char]10 + "Return the head and close this function call. "

# ---- Formatting: copy existing cell formats so the style table is reused (no duplicate styles) ----

# Column A category header fill (same style family as rows 1,2,4 / row 3)
$ws.Range("A1").Copy()
$ws.Range("A11").PasteSpecial(-4122)
$ws.Range("A1").Copy()
$ws.Range("A12").PasteSpecial(-4122)
$ws.Range("A3").Copy()
$ws.Range("A13").PasteSpecial(-4122)

# Column B style (used for B7:B10)
$ws.Range("B7").Copy()
$ws.Range("B11").PasteSpecial(-4122)
$ws.Range("B7").Copy()
$ws.Range("B12").PasteSpecial(-4122)
$ws.Range("B7").Copy()
$ws.Range("B13").PasteSpecial(-4122)

# Column C styles
$ws.Range("C2").Copy()
$ws.Range("C11").PasteSpecial(-4122)
$ws.Range("C3").Copy()
$ws.Range("C12").PasteSpecial(-4122)

$excel.CutCopyMode = 0

# C13 gets a new highlight fill (Blue, Accent 1, Lighter 60%) with top-aligned wrapped text
$ws.Range("C13").Interior.Color = 15189940
$ws.Range("C13").VerticalAlignment = -4160
$ws.Range("C13").WrapText = 1

# ---- Row heights ----
$ws.Rows.Item(11).RowHeight = 35.25
$ws.Rows.Item(12).RowHeight = 33.75
$ws.Rows.Item(13).RowHeight = 30

# ---- Selection ----
$ws.Range("A14").Select()
